# Fruta / hortaliza, semanal
# Insert two new weekly price records at the top of the "Durazno" data block
# (rows 206-243), shifting the existing 38 rows down to 208-245.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at row 206 - pushes rows 206:243 down to 208:245
$ws.Rows("206:207").Insert()

# Row 206 - new record: Carson / Especial
$ws.Cells.Item(206, 1).Value() = 4
$ws.Cells.Item(206, 2).Value() = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(206, 3).Value() = "Los Lagos"
$ws.Cells.Item(206, 4).Value() = 44946
$ws.Cells.Item(206, 5).Value() = 10
$ws.Cells.Item(206, 6).Value() = "Fruta"
$ws.Cells.Item(206, 7).Value() = 100103
$ws.Cells.Item(206, 8).Value() = "Frutos de hueso (carozo)"
$ws.Cells.Item(206, 9).Value() = 100103004
$ws.Cells.Item(206, 10).Value() = "Durazno"
$ws.Cells.Item(206, 11).Value() = "Carson"
$ws.Cells.Item(206, 12).Value() = "Especial"
$ws.Cells.Item(206, 13).Value() = 200
$ws.Cells.Item(206, 14).Value() = 22000
$ws.Cells.Item(206, 15).Value() = 22000
$ws.Cells.Item(206, 16).Value() = 22000
$ws.Cells.Item(206, 17).Value() = "$/caja 14 kilos empedrada"
$ws.Cells.Item(206, 18).Value() = "Región de O'Higgins"
$ws.Cells.Item(206, 19).Value() = 1571
$ws.Cells.Item(206, 20).Value() = 14

# Row 207 - new record: Carson / Primera
$ws.Cells.Item(207, 1).Value() = 4
$ws.Cells.Item(207, 2).Value() = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(207, 3).Value() = "Los Lagos"
$ws.Cells.Item(207, 4).Value() = 44946
$ws.Cells.Item(207, 5).Value() = 10
$ws.Cells.Item(207, 6).Value() = "Fruta"
$ws.Cells.Item(207, 7).Value() = 100103
$ws.Cells.Item(207, 8).Value() = "Frutos de hueso (carozo)"
$ws.Cells.Item(207, 9).Value() = 100103004
$ws.Cells.Item(207, 10).Value() = "Durazno"
$ws.Cells.Item(207, 11).Value() = "Carson"
$ws.Cells.Item(207, 12).Value() = "Primera"
$ws.Cells.Item(207, 13).Value() = 400
$ws.Cells.Item(207, 14).Value() = 18000
$ws.Cells.Item(207, 15).Value() = 19000
$ws.Cells.Item(207, 16).Value() = 18500
$ws.Cells.Item(207, 17).Value() = "$/caja 14 kilos empedrada"
$ws.Cells.Item(207, 18).Value() = "Región de O'Higgins"
$ws.Cells.Item(207, 19).Value() = 1321
$ws.Cells.Item(207, 20).Value() = 14
